$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F (想去人数 / interested count) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F17").Value = 170
$wsExhibit.Range("F18").Value = 215
$wsExhibit.Range("F22").Value = 2133
$wsExhibit.Range("F23").Value = 176
$wsExhibit.Range("F29").Value = 1201
$wsExhibit.Range("F30").Value = 4430
$wsExhibit.Range("F32").Value = 4056
$wsExhibit.Range("F33").Value = 1085
$wsExhibit.Range("F35").Value = 3130
$wsExhibit.Range("F37").Value = 1447
$wsExhibit.Range("F38").Value = 228
$wsExhibit.Range("F42").Value = 720
$wsExhibit.Range("F43").Value = 1106
$wsExhibit.Range("F46").Value = 451

# Sheet "全部类型" (all types) - same events, different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 170
$wsAll.Range("F15").Value = 215
$wsAll.Range("F18").Value = 2133
$wsAll.Range("F19").Value = 176
$wsAll.Range("F25").Value = 1201
$wsAll.Range("F28").Value = 4430
$wsAll.Range("F30").Value = 4056
$wsAll.Range("F31").Value = 1085
$wsAll.Range("F33").Value = 3130
$wsAll.Range("F36").Value = 1447
$wsAll.Range("F38").Value = 228
$wsAll.Range("F42").Value = 720
$wsAll.Range("F44").Value = 1106
$wsAll.Range("F47").Value = 451

$wb.Save()
